# Auto-generated Excel COM-interop script to apply the Zeromus_Profits update
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across several
# worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect refreshed market data.

$wb = $excel.ActiveWorkbook

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 348.7143
$ws.Range("I33").Value = 360.15384
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 360.15384
$ws.Range("L33").Value = 200
$ws.Range("M33").Value = -131.15384
$ws.Range("N33").Value = -658

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 304.9565
$ws.Range("I107").Value = 204.9375
$ws.Range("J107").Value = 533.5714
$ws.Range("K107").Value = 204.9375
$ws.Range("L107").Value = 533.5714
$ws.Range("M107").Value = 1715.0625
$ws.Range("N107").Value = -4373.5714

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1090265.1
$ws.Range("J112").Value = 1425616
$ws.Range("L112").Value = 4276848
$ws.Range("N112").Value = -4279064

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2876424.2
$ws.Range("I116").Value = 3207896.2
$ws.Range("J116").Value = 3666.6667
$ws.Range("K116").Value = 3207896.2
$ws.Range("L116").Value = 3666.6667
$ws.Range("M116").Value = -3204454.2
$ws.Range("N116").Value = -10550.6667

# ALC row 118
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 3997.8462
$ws.Range("I118").Value = 674.6667
$ws.Range("J118").Value = 7874.8887
$ws.Range("K118").Value = 2024.0001
$ws.Range("L118").Value = 23624.6661
$ws.Range("M118").Value = -367.0001
$ws.Range("N118").Value = -26938.6661

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3538.5305
$ws.Range("I138").Value = 1966.3939
$ws.Range("J138").Value = 4336.6924
$ws.Range("K138").Value = 5899.1817
$ws.Range("L138").Value = 13010.0772
$ws.Range("M138").Value = -759.1817000000001
$ws.Range("N138").Value = -23290.0772

# ARM row 24
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 24150
$ws.Range("J24").Value = 24150
$ws.Range("L24").Value = 24150
$ws.Range("N24").Value = -24898

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2934.1667
$ws.Range("I45").Value = 3368.5
$ws.Range("J45").Value = 2499.8333
$ws.Range("K45").Value = 3368.5
$ws.Range("L45").Value = 2499.8333
$ws.Range("M45").Value = -2991.5
$ws.Range("N45").Value = -3253.8333

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2578.6
$ws.Range("I74").Value = 2587.5173
$ws.Range("J74").Value = 2535.5
$ws.Range("K74").Value = 2587.5173
$ws.Range("L74").Value = 2535.5
$ws.Range("M74").Value = -1713.5173
$ws.Range("N74").Value = -4283.5

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2578.6
$ws.Range("I77").Value = 2587.5173
$ws.Range("J77").Value = 2535.5
$ws.Range("K77").Value = 12937.5865
$ws.Range("L77").Value = 12677.5
$ws.Range("M77").Value = -8569.586499999999
$ws.Range("N77").Value = -21413.5

# ARM row 100
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H100").Value = 24150
$ws.Range("J100").Value = 24150
$ws.Range("L100").Value = 24150
$ws.Range("N100").Value = -26314

# ARM row 103
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H103").Value = 32000
$ws.Range("J103").Value = 32000
$ws.Range("L103").Value = 32000
$ws.Range("N103").Value = -34344

# ARM row 108
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H108").Value = 52500
$ws.Range("J108").Value = 52500
$ws.Range("L108").Value = 52500
$ws.Range("N108").Value = -60180

# ARM row 109
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 29800
$ws.Range("J109").Value = 29800
$ws.Range("L109").Value = 29800
$ws.Range("N109").Value = -32574

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1546.7
$ws.Range("I110").Value = 1544.909
$ws.Range("J110").Value = 1551.625
$ws.Range("K110").Value = 1544.909
$ws.Range("L110").Value = 1551.625
$ws.Range("M110").Value = 500.0909999999999
$ws.Range("N110").Value = -5641.625

# ARM row 112
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 29693.5
$ws.Range("J112").Value = 29693.5
$ws.Range("L112").Value = 29693.5
$ws.Range("N112").Value = -32647.5

# ARM row 115
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").ClearContents()
$ws.Range("N115").ClearContents()

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1824.7727
$ws.Range("I105").Value = 1403.3334
$ws.Range("J105").Value = 2727.8572
$ws.Range("K105").Value = 1403.3334
$ws.Range("L105").Value = 2727.8572
$ws.Range("M105").Value = 343.6666
$ws.Range("N105").Value = -6221.8572

# BSM row 110
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 42000
$ws.Range("J110").Value = 42000
$ws.Range("L110").Value = 42000
$ws.Range("N110").Value = -50180

# BSM row 135
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 46575
$ws.Range("J135").Value = 46575
$ws.Range("L135").Value = 46575
$ws.Range("N135").Value = -56715

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 66.181816
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 66.181816
$ws.Range("K7").Value = 0
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -292.181816

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1066.3
$ws.Range("I94").Value = 947.8
$ws.Range("J94").Value = 1184.8
$ws.Range("K94").Value = 947.8
$ws.Range("L94").Value = 1184.8
$ws.Range("M94").Value = -496.8
$ws.Range("N94").Value = -2086.8

# CUL row 32
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 2500
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

# CUL row 120
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 13200
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 13200
$ws.Range("K120").Value = 0
$ws.Range("L120").ClearContents()
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -49276

# GSM row 108
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").ClearContents()
$ws.Range("N108").ClearContents()

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2133.3333
$ws.Range("I7").Value = 2150
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 2150
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -2038
$ws.Range("N7").Value = -2224

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 609.85
$ws.Range("I16").Value = 557.46155
$ws.Range("J16").Value = 707.1429000000001
$ws.Range("K16").Value = 557.46155
$ws.Range("L16").Value = 707.1429000000001
$ws.Range("M16").Value = -387.46155
$ws.Range("N16").Value = -1047.1429

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1582.5454
$ws.Range("I93").Value = 1733
$ws.Range("J93").Value = 1432.091
$ws.Range("K93").Value = 1733
$ws.Range("L93").Value = 1432.091
$ws.Range("M93").Value = -485
$ws.Range("N93").Value = -3928.091

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2133.3333
$ws.Range("I126").Value = 2150
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 6450
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -3980
$ws.Range("N126").Value = -10940

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1043.5294
$ws.Range("I107").Value = 1226.6923
$ws.Range("J107").Value = 448.25
$ws.Range("K107").Value = 3680.0769
$ws.Range("L107").Value = 1344.75
$ws.Range("M107").Value = -1760.0769
$ws.Range("N107").Value = -5184.75
